# This edit inserts a new price-report record (row) into the weekly
# "Apio" (celery) price sheet for "Macroferia Regional de Talca", at row 43,
# shifting all subsequent records down by one row (old row 43 becomes new
# row 44, old row 192 becomes new row 193, etc.), and populates the new
# row 43 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 43; this pushes rows 43-192 down to 44-193
# (Excel automatically extends the used range / dimension to R193).
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with the new record's data.
$ws.Cells.Item(43, 1).Value  = 5
$ws.Cells.Item(43, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(43, 3).Value  = "Maule"
$ws.Cells.Item(43, 4).Value  = 44742
$ws.Cells.Item(43, 5).Value  = 7
$ws.Cells.Item(43, 6).Value  = 100112017
$ws.Cells.Item(43, 7).Value  = "Apio"
$ws.Cells.Item(43, 8).Value  = "Americana (o)"
$ws.Cells.Item(43, 9).Value  = "Primera"
$ws.Cells.Item(43, 10).Value = 600
$ws.Cells.Item(43, 11).Value = 7000
$ws.Cells.Item(43, 12).Value = 7000
$ws.Cells.Item(43, 13).Value = 7000
$ws.Cells.Item(43, 14).Value = "`$/docena de matas"
$ws.Cells.Item(43, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(43, 16).Value = 1167
$ws.Cells.Item(43, 17).Value = 6
$ws.Cells.Item(43, 18).Value = "Hortaliza"
